$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 186, shifting the existing rows 186..236 down to 187..237.
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new weekly price observation.
$ws.Cells.Item(186, 1).Value = 4
$ws.Cells.Item(186, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value = "Los Lagos"
$ws.Cells.Item(186, 4).Value = 44754
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = 100112039
$ws.Cells.Item(186, 7).Value = "Ciboulette"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 180
$ws.Cells.Item(186, 11).Value = 3500
$ws.Cells.Item(186, 12).Value = 3500
$ws.Cells.Item(186, 13).Value = 3500
$ws.Cells.Item(186, 14).Value = "$/docena de atados"
$ws.Cells.Item(186, 15).Value = "Región Metropolitana"
$ws.Cells.Item(186, 16).Value = 1167
$ws.Cells.Item(186, 17).Value = 3
$ws.Cells.Item(186, 18).Value = "Hortaliza"
